{"js": "// Renumber the \"AC ID\" column of the Acceptance Criteria table.\n// Each old value (format X.1.Y) is replaced by the new value (format X.Y),\n// i.e. the middle \".1\" segment is dropped.\nconst idMap = [\n  [\"1.1.1\", \"1.1\"],\n  [\"1.1.2\", \"1.2\"],\n  [\"1.1.3\", \"1.3\"],\n  [\"1.1.4\", \"1.4\"],\n  [\"1.1.5\", \"1.5\"],\n  [\"2.1.1\", \"2.1\"],\n  [\"2.1.2\", \"2.2\"],\n  [\"3.1.0\", \"3.0\"],\n  [\"4.1.1\", \"4.1\"],\n  [\"4.1.2\", \"4.2\"],\n  [\"4.1.3\", \"4.3\"],\n  [\"4.1.4\", \"4.4\"],\n  [\"5.1.1\", \"5.1\"],\n  [\"5.1.2\", \"5.2\"],\n  [\"5.1.3\", \"5.3\"],\n  [\"6.1.1\", \"6.1\"],\n  [\"6.1.2\", \"6.2\"],\n  [\"6.1.3\", \"6.3\"],\n  [\"6.1.4\", \"6.4\"],\n  [\"6.1.5\", \"6.5\"],\n  [\"6.1.6\", \"6.6\"],\n  [\"7.1.1\", \"7.1\"],\n  [\"7.1.2\", \"7.2\"],\n];\n\nconst body = context.document.body;\n\n// First, locate all the search results (one per old value) so we only load\n// what we need.\nconst searchResults = idMap.map(([oldVal]) =>\n  body.search(oldVal, { matchCase: true, matchWholeWord: true })\n);\nsearchResults.forEach((r) => r.load(\"items/text\"));\nawait context.sync();\n\nfor (let i = 0; i < idMap.length; i++) {\n  const [, newVal] = idMap[i];\n  const results = searchResults[i];\n  for (const range of results.items) {\n    range.insertText(newVal, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Renumber the \"AC ID\" column of the Acceptance Criteria table.\n# Each old value (format X.1.Y) is replaced by the new value (format X.Y),\n# i.e. the middle \".1\" segment is dropped.\n$d = $word.ActiveDocument\n\n$idMap = @(\n    @(\"1.1.1\", \"1.1\"),\n    @(\"1.1.2\", \"1.2\"),\n    @(\"1.1.3\", \"1.3\"),\n    @(\"1.1.4\", \"1.4\"),\n    @(\"1.1.5\", \"1.5\"),\n    @(\"2.1.1\", \"2.1\"),\n    @(\"2.1.2\", \"2.2\"),\n    @(\"3.1.0\", \"3.0\"),\n    @(\"4.1.1\", \"4.1\"),\n    @(\"4.1.2\", \"4.2\"),\n    @(\"4.1.3\", \"4.3\"),\n    @(\"4.1.4\", \"4.4\"),\n    @(\"5.1.1\", \"5.1\"),\n    @(\"5.1.2\", \"5.2\"),\n    @(\"5.1.3\", \"5.3\"),\n    @(\"6.1.1\", \"6.1\"),\n    @(\"6.1.2\", \"6.2\"),\n    @(\"6.1.3\", \"6.3\"),\n    @(\"6.1.4\", \"6.4\"),\n    @(\"6.1.5\", \"6.5\"),\n    @(\"6.1.6\", \"6.6\"),\n    @(\"7.1.1\", \"7.1\"),\n    @(\"7.1.2\", \"7.2\")\n)\n\nforeach ($pair in $idMap) {\n    $oldVal = $pair[0]\n    $newVal = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Text = $oldVal\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $true\n    $find.MatchWildcards = $false\n\n    $found = $find.Execute()\n    if ($found) {\n        $range.Text = $newVal\n    } else {\n        Write-Output \"WARNING: could not find AC ID '$oldVal'\"\n    }\n}\n"}
